$d = $word.ActiveDocument

# --- Change 1: remove the stray _GoBack bookmark that wraps "organizers" ---
$bm = $d.Bookmarks("_GoBack")
$bm.Delete()

# --- Change 2: rewrite the "Please note..." paragraph and split off the
#     trailing page-break run into its own paragraph, re-adding a fresh
#     _GoBack bookmark in the middle of the rewritten text. ---

$found = $d.Content.Find.Execute("Please note that if you can")
if (-not $found) {
    throw "Could not locate target paragraph"
}
$hitStart = $d.Content.Start

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Please note that if you can*") {
        $targetPara = $cand
        break
    }
}
if ($targetPara -eq $null) {
    throw "Could not find paragraph object"
}

$rng = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

$ooxml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
'<w:p><w:pPr><w:jc w:val="both"/><w:rPr><w:b/><w:color w:val="FF0000"/><w:sz w:val="24"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:b/><w:color w:val="FF0000"/><w:sz w:val="24"/></w:rPr><w:t>Please note that if you can use only a non-Hungarian bank for transfer, then you are not required to pay deposit - due to high amount of tra</w:t></w:r>' +
'<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
'<w:r><w:rPr><w:b/><w:color w:val="FF0000"/><w:sz w:val="24"/></w:rPr><w:t>nsfer costs. You can settle your registration fee at the conference registration venue in cash. In this case, please expect an e-mail from the organizers right before the conference when we ask you for a second confirmation of your application.</w:t></w:r>' +
'</w:p>' +
'<w:p><w:pPr><w:spacing w:before="120" w:after="120"/><w:jc w:val="both"/><w:rPr><w:color w:val="FF0000"/><w:sz w:val="28"/></w:rPr></w:pPr>' +
'<w:r><w:rPr><w:color w:val="FF0000"/><w:sz w:val="28"/></w:rPr><w:br w:type="page"/></w:r></w:p>' +
'</w:body></w:document>' +
'</pkg:xmlData></pkg:part></pkg:package>'

$rng.InsertXML($ooxml)

Write-Output "done"
